$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.963.70'
$ws.Range("E2").Value = '  +0.99%  '
$ws.Range("D3").Value = '3.253.00'
$ws.Range("E3").Value = '  +1.70%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range("E4").Value = '  -0.15%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '606.86'
$c.ClearFormats()
$ws.Range("E5").Value = '  +1.69%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '157.69'
$c.ClearFormats()
$ws.Range("E6").Value = '  +2.25%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.254.09'
$ws.Range("E8").Value = '  +1.86%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.548'
$c.ClearFormats()
$ws.Range("E9").Value = '  +2.26%  '
$ws.Range("E10").Value = '  +0.37%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.70'
$c.ClearFormats()
$ws.Range("E11").Value = '  -6.70%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.513'
$c.ClearFormats()
$ws.Range("E12").Value = '  -0.13%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000273'
$c.ClearFormats()
$ws.Range("E13").Value = '  +0.94%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '39.20'
$c.ClearFormats()
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").Value = '3.788.42'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").Value = '66.876.56'
$ws.Range("E16").Value = '  +0.93%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '7.46'
$c.ClearFormats()
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '3.246.30'
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("E19").Value = '  +1.12%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '514.31'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.78%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '15.44'
$c.ClearFormats()
$ws.Range("E21").Value = '  +0.62%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.741'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.48%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '8.13'
$c.ClearFormats()
$ws.Range("E23").Value = '  +1.54%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '14.93'
$c.ClearFormats()
$ws.Range("E24").Value = '  -1.17%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '85.04'
$c.ClearFormats()
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  +0.18%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.47'
$c.ClearFormats()
$ws.Range("E27").Value = '  +1.82%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.ClearFormats()
$ws.Range("E28").Value = '  +0.93%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.ClearFormats()
$ws.Range("E29").Value = '  +6.03%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '3.06'
$c.ClearFormats()
$ws.Range("E30").Value = '  +5.08%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.15'
$c.ClearFormats()
$ws.Range("E31").Value = '  +2.67%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '28.36'
$c.ClearFormats()
$ws.Range("E32").Value = '  +0.43%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range("E33").Value = '  -0.15%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.19'
$c.ClearFormats()
$ws.Range("E34").Value = '  -2.73%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '6.56'
$c.ClearFormats()
$ws.Range("E35").Value = '  +0.45%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '528.08'
$c.ClearFormats()
$ws.Range("E36").Value = '  +9.03%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '56.29'
$c.ClearFormats()
$ws.Range("E37").Value = '  +2.75%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0934'
$c.ClearFormats()
$ws.Range("E38").Value = '  +3.77%  '
$ws.Range("D39").Value = '0.0₃0770'
$ws.Range("E39").Value = '  +17.81%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0422'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.01'
$c.ClearFormats()
$ws.Range("E41").Value = '  +3.45%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.128'
$c.ClearFormats()
$ws.Range("E42").Value = '  +5.08%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '8.85'
$c.ClearFormats()
$ws.Range("E43").Value = '  +0.11%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.304'
$c.ClearFormats()
$ws.Range("E44").Value = '  +1.65%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.ClearFormats()
$ws.Range("E45").Value = '  +4.87%  '
$ws.Range("D46").Value = '2.872.77'
$ws.Range("E46").Value = '  -1.78%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '28.72'
$c.ClearFormats()
$ws.Range("E47").Value = '  +0.91%  '
$ws.Range("E48").Value = '  +3.99%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("E51").Value = '  +0.72%  '
